$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") date value moves from 45205 (2023-10-06) to
# 45206 (2023-10-07) for every data row (rows 2 through 224).
$ws.Range("C2:C224").Value = 45206
